# Apply the "Add files via upload" revision to the CityEats ALS deck.
# Strategy: use TextRange.Find()+assignment for in-place text swaps so we
# don't disturb paragraphs/runs we're not touching, and TextRange.InsertBefore("`r")
# to add a leading blank paragraph cleanly (matches the shape PowerPoint itself
# produces for a brand-new leading empty bullet line).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - Title slide: update the subtitle/byline text.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "MET CS 777 - Srivatsav Shrikanth (sri99-svg)"

# ---------------------------------------------------------------------------
# Slide 2 - Agenda: three of the five bullet lines are reworded.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2body = $s2.Shapes.Item(2)

$f = $s2body.TextFrame.TextRange.Find("How to Run")
$f.Text = "Pipeline & Repro"

$f = $s2body.TextFrame.TextRange.Find("Results (K=50)")
$f.Text = "Results (P@K / R@K / NDCG@K)"

$f = $s2body.TextFrame.TextRange.Find("Next Steps")
$f.Text = "Conclusion & Next Steps"

# ---------------------------------------------------------------------------
# Slide 3 - was "How to Run (quickcheck)", now "Problem & Data".
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$f = $s3.Shapes.Item(1).TextFrame.TextRange.Find("How to Run (quickcheck)")
$f.Text = "Problem & Data"

$s3.Shapes.Item(2).TextFrame.TextRange.Text = "Goal: top-N recommendations from explicit ratings.`rRatings schema: user_id, item_id, rating`rLarge scale (multi-million rows)`rParquet input; Spark for scale"

# ---------------------------------------------------------------------------
# Slide 4 - was "Results (K=50)", now "Methodology (ALS)". Content is fully
# different (4 bullets -> blank + 3 bullets) so rebuild the body outright.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4title = $s4.Shapes.Item(1)
$s4body = $s4.Shapes.Item(2)

$f = $s4title.TextFrame.TextRange.Find("Results (K=50)")
$f.Text = "Methodology (ALS)"

$s4body.TextFrame.TextRange.Text = "`rStringIndexer -> (user_idx, biz_idx)`rALS rank/regParam/maxIter; nonnegative, coldStart=drop`rSplit: global random (seed 42)"

# ---------------------------------------------------------------------------
# Slide 5 - was "Next Steps", now "Pipeline & Repro".
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5title = $s5.Shapes.Item(1)
$s5body = $s5.Shapes.Item(2)

$f = $s5title.TextFrame.TextRange.Find("Next Steps")
$f.Text = "Pipeline & Repro"

$f = $s5body.TextFrame.TextRange.Find("Run on Dataproc Serverless")
$f.Text = "jobs/train_als_local.py + conf/config.yaml"

$f = $s5body.TextFrame.TextRange.Find("Min-interaction filters (>=5)")
$f.Text = "spark-submit command in README"

$f = $s5body.TextFrame.TextRange.Find("Light rank/reg sweep")
$f.Text = "Spark UI at localhost:4040"

$s5body.TextFrame.TextRange.InsertBefore("`r") | Out-Null

# ---------------------------------------------------------------------------
# New slide 6 - "Results".
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Results"
$s6body = $s6.Shapes.Item(2)
$s6body.TextFrame.TextRange.Text = "K=50 (quickcheck)`rPrecision@K: <value>`rRecall@K: <value>`rNDCG@K: <value>"
$s6body.TextFrame.TextRange.InsertBefore("`r") | Out-Null

# ---------------------------------------------------------------------------
# New slide 7 - "Artifacts".
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Artifacts"
$s7body = $s7.Shapes.Item(2)
$s7body.TextFrame.TextRange.Text = "Model, UI/Item maps`rmetrics.json, manifest.json`rBatch CSV / Recs JSON (optional)"
$s7body.TextFrame.TextRange.InsertBefore("`r") | Out-Null

# ---------------------------------------------------------------------------
# New slide 8 - "Conclusion & Next Steps".
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusion & Next Steps"
$s8body = $s8.Shapes.Item(2)
$s8body.TextFrame.TextRange.Text = "Good ranking quality on explicit ratings`rNext: hyperparameter sweep & cluster run`rServing API + monitoring"
$s8body.TextFrame.TextRange.InsertBefore("`r") | Out-Null
